$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header text: "unnamed: 1_level_1" -> "total"
$ws.Range("B2").Value = "total"

# Remove the now-unneeded rows, working from the bottom up so the row
# numbers of earlier rows are not shifted by the later deletions:
#   row 41 -> footnote "fonte: ibge, diretoria de pesquisas, ..."
#   row 8  -> subheading "grandes regioes e unidades da federacao"
#   row 5  -> subheading "situacao do domicilio"
$ws.Rows.Item(41).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
